$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13: Geminox 320mg Tablet - 8's -------------------------------
$ws.Range("H13").Value = 72

# --- UOM (column D) corrections for Ketonic / Kynol rows --------------
$ws.Range("D15").Value = "5 's"     # Ketonic 10mg Tablet
$ws.Range("D16").Value = "20's"     # Ketonic 30mg IM/IV Injection - 4's
$ws.Range("D18").Value = "50 's"    # Kynol TR 200mg Capsule
$ws.Range("D19").Value = "30 's"    # Kynol TR 100mg Capsule

# --- Row 23: Sk-Mox 500mg Capsule stock/sales figures ------------------
$ws.Range("E23").Value = 2
$ws.Range("I23").Value = 7
$ws.Range("L23").Value = 25
$ws.Range("N23").Value = 175
$ws.Range("O23").Value = 182
$ws.Range("T23").Value = 182
$ws.Range("AL23").Value = 25
$ws.Range("AZ23").Value = 46
$ws.Range("BA23").Value = 10079
